$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (CovidDeclaration), shifting
# CovidDeclaration and every later column one position to the right.
$ws.Range("G1:G2").EntireColumn.Insert()

# Populate the new column G with the "Signed" header/value pair.
$ws.Range("G1").Value = "Signed"
$ws.Range("G2").Value = "&=FundingClaimDataExtractInfo.FundingClaimsDataExtract.Signed"

# The new column keeps the same custom width as its neighbouring columns.
$ws.Range("G1").ColumnWidth = 9

# Refresh the sheet's remembered sort state (Data > Sort) so the sort
# key columns that were pushed right by the insert (old H -> I, old J -> K)
# point at their new locations, and the remembered range now also spans
# the newly inserted column.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D2:D266")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("C2:C266")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("A2:A266")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("I2:I266")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("K2:K266")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:V266"))
$ws.Sort.Header = -4163
$ws.Sort.Apply()

# Update the selection to match the post-edit workbook state.
$ws.Range("F4").Select()
